# Adding labs 13 and 11: append a second "HumMod"/"QCP" data block
# (rows 11-20) below the existing "Salt Variance" table on Sheet1,
# mirroring the layout of rows 1-10, plus a new "QCP" label at D1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section title row (row 11), styled like the other data-row
#     labels (Arial font, wrap text, vertical-top) but without a box
#     border - matches font/alignment applied by the source rows.
#     ("HumMod" must become shared-string index 10, added before
#     "QCP" below, to match the target string table order.) ---
$ws.Range("A3").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "HumMod"
$ws.Range("A11").Borders.LineStyle = -4142
$ws.Rows("11").RowHeight = 30.75

# --- New header cell on row 1 ("QCP" becomes shared-string index 11) ---
$ws.Range("D1").Value = "QCP"

# --- Row 12: mirrors row 2 (Sodium Intake header) ---
$ws.Range("A2:D2").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("F2:I2").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Rows("12").RowHeight = 30.75

$ws.Range("A12").Value = "Sodium Intake"
$ws.Range("B12").Value = 20
$ws.Range("C12").Value = 180
$ws.Range("D12").Value = 500
$ws.Range("F12").Value = "Sodium Intake"
$ws.Range("G12").Value = 20
$ws.Range("H12").Value = 180
$ws.Range("I12").Value = 500

# --- Row 13: mirrors row 3 (Arterial Pressure) ---
$ws.Range("A3:D3").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("F3:I3").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Rows("13").RowHeight = 60.75

$ws.Range("A13").Value = "Arterial Pressure(mmHg)"
$ws.Range("B13").Value = 97.3
$ws.Range("C13").Value = 96.6
$ws.Range("D13").Value = 96.9
$ws.Range("F13").Value = "Arterial Pressure(mmHg)"
$ws.Range("G13").Value = 84.8
$ws.Range("H13").Value = 89
$ws.Range("I13").Value = 90.9

# --- Row 14: mirrors row 4 (Right Atrial Pressure) ---
$ws.Range("A4:D4").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("F4:I4").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Rows("14").RowHeight = 75.75

$ws.Range("A14").Value = "Right Atrial Pressure(mmHg)"
$ws.Range("B14").Value = 0.1
$ws.Range("C14").Value = 1.2
$ws.Range("D14").Value = 2.1
$ws.Range("F14").Value = "Right Atrial Pressure(mmHg)"
$ws.Range("G14").Value = -1
$ws.Range("H14").Value = 0.6
$ws.Range("I14").Value = 1.6

# --- Row 15: mirrors row 5 (Left Atrial Pressure) ---
$ws.Range("A5:D5").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("F5:I5").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Rows("15").RowHeight = 75.75

$ws.Range("A15").Value = "Left Atrial Pressure(mmHg)"
$ws.Range("B15").Value = 4.1
$ws.Range("C15").Value = 5.2
$ws.Range("D15").Value = 6.1
$ws.Range("F15").Value = "Left Atrial Pressure(mmHg)"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 5.1

# --- Row 16: mirrors row 6 (Plasma [AngII]) ---
$ws.Range("A6:D6").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("F6:I6").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Rows("16").RowHeight = 45.75

$ws.Range("A16").Value = "Plasma [AngII](pg/mL)"
$ws.Range("B16").Value = 40.4
$ws.Range("C16").Value = 22.3
$ws.Range("D16").Value = 14.2
$ws.Range("F16").Value = "Plasma [AngII](pg/mL)"
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0

# --- Row 17: mirrors row 7 (Plasma [Aldosterone]) ---
$ws.Range("A7:D7").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("F7:I7").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Rows("17").RowHeight = 60.75

$ws.Range("A17").Value = "Plasma [Aldosterone](pmol/L)"
$ws.Range("B17").Value = 519.7
$ws.Range("C17").Value = 306.9
$ws.Range("D17").Value = 215.1
$ws.Range("F17").Value = "Plasma [Aldosterone](pmol/L)"
$ws.Range("G17").Value = 249.8
$ws.Range("H17").Value = 179.3
$ws.Range("I17").Value = 146

# --- Row 18: mirrors row 8 (Plasma [ANP]) ---
$ws.Range("A8:D8").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("F8:I8").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Rows("18").RowHeight = 45.75

$ws.Range("A18").Value = "Plasma [ANP](pmol/L)"
$ws.Range("B18").Value = 19.6
$ws.Range("C18").Value = 22.4
$ws.Range("D18").Value = 24.8
$ws.Range("F18").Value = "Plasma [ANP](pmol/L)"
$ws.Range("G18").Value = 14.1
$ws.Range("H18").Value = 21.9
$ws.Range("I18").Value = 25

# --- Row 19: mirrors row 9 (Urine Na+ Excretion) - note A:D use the
#     same style as rows 13-18, while F:I use the "last data row"
#     style (like F9:I9) that's also merged with the row below. ---
$ws.Range("A9:D9").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("F9:I9").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Rows("19").RowHeight = 75.75

$ws.Range("A19").Value = "Urine Na+ Excretion(mEq/min)"
$ws.Range("B19").Value = 0.029
$ws.Range("C19").Value = 0.123
$ws.Range("D19").Value = 0.359
$ws.Range("F19").Value = "Urine Na+ Excretion(mEq/min)"
$ws.Range("G19").Value = 0.023
$ws.Range("H19").Value = 0.125
$ws.Range("I19").Value = 0.342

# --- Row 20: mirrors row 10 (closing thick-bottom border row) ---
$ws.Range("F10:I10").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Rows("20").RowHeight = 15.75

# --- Merge the second block's right-hand label column with its
#     closing row, same as F9:F10 / G9:G10 / H9:H10 / I9:I10. ---
$ws.Range("F19:F20").Merge()
$ws.Range("G19:G20").Merge()
$ws.Range("H19:H20").Merge()
$ws.Range("I19:I20").Merge()

# --- Selection / view state matches the post-edit workbook: the new
#     block is in view (no frozen/top-left scroll override) and the
#     freshly-pasted second table is selected. ---
$ws.Range("F12:I20").Select()
